{"js": "// Replace four English interview prompts with their Afrikaans translations.\nconst replacements = [\n  {\n    find: \"I want to go through some of the situations which your X might have spoken to you about when it comes to keeping safe from sexual violence. \",\n    replace: \"Ek wil deur 'n paar situasies gaan wat jou X dalk met jou bespreek het, oor om jou veilig te hou van seksuele geweld. \"\n  },\n  {\n    find: \"Have you had conversations about this since?\",\n    replace: \"Het jy gespreke gehad oor dit sedert dan?\"\n  },\n  {\n    find: \"Has your X ever talked to you about having to do anything sexual that you didn\\u2019t want to do?\",\n    replace: \"Het jou X ooit met jou gepraat oor om iets seksueel te doen wat jy nie wou nie?\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.Text = $findText\n    $find.MatchCase = $true\n    $find.MatchWildcards = $false\n    while ($find.Execute()) {\n        $rng.Text = $replaceText\n        $rng.Collapse(0)\n        $rng.End = $d.Content.End\n    }\n}\n\nReplace-AllText \"I want to go through some of the situations which your X might have spoken to you about when it comes to keeping safe from sexual violence. \" \"Ek wil deur 'n paar situasies gaan wat jou X dalk met jou bespreek het, oor om jou veilig te hou van seksuele geweld. \"\nReplace-AllText \"Have you had conversations about this since?\" \"Het jy gespreke gehad oor dit sedert dan?\"\nReplace-AllText \"Has your X ever talked to you about having to do anything sexual that you didn\u2019t want to do?\" \"Het jou X ooit met jou gepraat oor om iets seksueel te doen wat jy nie wou nie?\"\n"}
